$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.994.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.81"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.501"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.47"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.61%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.00%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.657.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.298.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.802"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.850.07"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.58"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.07"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.84"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.23"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.94"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0695"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.75"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.966.95"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.67"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.77"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.525.34"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.05%  "
